$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 16, shifting rows 16-26 down to 17-27,
# copying formatting from the row above (row 15) by using the default
# Insert behaviour (format follows the row above in Excel).
$ws.Rows.Item(16).Insert()

# Populate the new criterion row (Talc)
$ws.Cells.Item(16, 1).Value = "F-2080"
$ws.Cells.Item(16, 2).Value = "talc"
$ws.Cells.Item(16, 3).Value = 2070
$ws.Cells.Item(16, 4).Value = 2100
$ws.Cells.Item(16, 5).Value = "hull corrected"

# Match the formatting used by the rest of the criteria rows (A col style,
# B:E col style) by copying it down from the row that used to be row 16
# (now row 17) which already carries the correct look.
$ws.Range("A17:E17").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-set the values (PasteSpecial formats only shouldn't touch values, but
# make sure they are exactly as intended)
$ws.Cells.Item(16, 1).Value = "F-2080"
$ws.Cells.Item(16, 2).Value = "talc"
$ws.Cells.Item(16, 3).Value = 2070
$ws.Cells.Item(16, 4).Value = 2100
$ws.Cells.Item(16, 5).Value = "hull corrected"

# Update selection to match the authored workbook state
$ws.Range("A16").Select()
